$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new BOM rows (5V charge-pump motor-control parts) ---
# Row 12: new decoupling-cap line (C17/C18/C19, 1uF)
$ws.Rows.Item(12).Insert()
# Row 33 (after the first insert shifts everything down by one): new CAT3200-5 IC line
$ws.Rows.Item(33).Insert()

# --- Populate new row 12 ---
$ws.Range("A12").Value = "C18, C19, C17, "
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "1uF"
$ws.Range("D12").Value = "Capacitors_SMD:C_0603"
$ws.Range("E12").Value = "Unpolarized capacitor"
$ws.Range("F12").Value = "311-3484-1-ND"
$ws.Range("G12").Value = 0.0582
$ws.Range("H12").Formula = "=G12*B12"

# --- Populate new row 33 ---
$ws.Range("A33").Value = "U5, "
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "CAT3200-5"
$ws.Range("D33").Value = "TO_SOT_Packages_SMD:SOT-23-6"
$ws.Range("E33").Value = "Charge Pump Switching Regulator IC Positive Fixed 5V 1 Output 100mA SOT-23-6"
$ws.Range("F33").Value = "CAT3200TDI-GT3OSCT-ND"
$ws.Range("G33").Value = 1.367
$ws.Range("H33").Formula = "=G33*B33"

# --- Update the "Generated:" timestamp banner ---
$ws.Range("A2").Value = "Generated: 1/6/2018  3:56:19 PM"

# --- Update sheet view / selection (was topLeftCell B4 / selection C18) ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()

$wb.Save()
